# New Orleans shard update:
#  1. hotel_info gains a new "State" column (inserted between Hotel_Name and
#     City) populated with "Louisiana" for the existing data row.
#  2. The sheet tab order is swapped so "review_info" appears before
#     "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the State column into hotel_info -----------------------
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder sheets: review_info first, hotel_info second ----------
$review = $wb.Worksheets.Item("review_info")
$review.Move($wb.Worksheets.Item(1))

Write-Output "done"
